$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.405.81'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.453.62'
$ws.Range("E3").Value = '  +1.29%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.51'
$ws.Range("E5").Value = '  +2.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.55'
$ws.Range("E6").Value = '  +0.30%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +0.77%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.449.28'
$ws.Range("E9").Value = '  +1.20%  '
$ws.Range("E10").Value = '  +0.73%  '
$ws.Range("E11").Value = '  +2.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.21'
$ws.Range("E12").Value = '  +0.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.347'
$ws.Range("E13").Value = '  -1.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.43'
$ws.Range("E14").Value = '  -0.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000177'
$ws.Range("E15").Value = '  +2.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.897.60'
$ws.Range("E16").Value = '  +1.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.214.17'
$ws.Range("E17").Value = '  -0.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.448.54'
$ws.Range("E18").Value = '  +0.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.88'
$ws.Range("E19").Value = '  -2.14%  '
$ws.Range("E20").Value = '  -0.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '328.67'
$ws.Range("E21").Value = '  +1.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.14'
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.97'
$ws.Range("E23").Value = '  -4.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.65'
$ws.Range("E25").Value = '  +0.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.31'
$ws.Range("E26").Value = '  +4.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '590.65'
$ws.Range("E27").Value = '  -3.55%  '
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0967'
$ws.Range("E28").Value = '  -0.95%  '
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.576.12'
$ws.Range("E29").Value = '  +1.03%  '
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.44'
$ws.Range("E31").Value = '  -1.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.05'
$ws.Range("E32").Value = '  -0.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.89'
$ws.Range("E33").Value = '  +0.92%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.94'
$ws.Range("E35").Value = '  -1.98%  '
$ws.Range("E36").Value = '  +0.35%  '
$ws.Range("E37").Value = '  -1.45%  '
$ws.Range("E38").Value = '  +1.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '154.42'
$ws.Range("E39").Value = '  +4.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.34'
$ws.Range("E40").Value = '  +1.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.43'
$ws.Range("E41").Value = '  -1.45%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.73'
$ws.Range("E42").Value = '  -0.46%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.05'
$ws.Range("E43").Value = '  +2.24%  '
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.46'
$ws.Range("E45").Value = '  -2.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '143.03'
$ws.Range("E46").Value = '  -1.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.66'
$ws.Range("E47").Value = '  -1.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0255'
$ws.Range("E48").Value = '  +17.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.608'
$ws.Range("E49").Value = '  +2.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0523'
$ws.Range("E50").Value = '  -0.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.85'
$ws.Range("E51").Value = '  -2.25%  '
